# Append three new diary entries (12/15/16 joulu) to the bottom of the log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 44: 12 joulu ----------------------------------------------------
$ws.Range("A44").Value = "12 joulu"

$ws.Range("B44").NumberFormat = "h:mm"
$ws.Range("B44").Value = "15.00-16.00"

$ws.Range("C44").WrapText = $true
$ws.Range("C44").Value = "Viimeinen tsemppirinki"

$ws.Range("G44").Value = 1

# ---- Row 45: 15 joulu ------------------------------------------------------
$ws.Range("A45").Value = "15 joulu"

$ws.Range("B45").NumberFormat = "h:mm"
$ws.Range("B45").Value = "18.30-20.00"

$ws.Range("C45").WrapText = $true
$ws.Range("C45").Value = "Tutkimusartikkelin läpisilmäily, "

$ws.Range("D45").WrapText = $true
$ws.Range("D45").Value = "Vähän tällaistahan tämä on, oppii kyllä hyvin uutta mutta tähän olisi pitänyt varata enemmän aikaa."

$ws.Range("G45").Value = 1.5

$ws.Rows.Item(45).RowHeight = 45

# ---- Row 46: 16 joulu ------------------------------------------------------
$ws.Range("A46").Value = "16 joulu"

$ws.Range("B46").NumberFormat = "h:mm"
$ws.Range("B46").WrapText = $true
$ws.Range("B46").Value = "10.00-11.15, 13.00-13.45, 17.30-18.30"

$ws.Range("C46").WrapText = $true
$ws.Range("C46").Value = "Tutkimusartikkelin taustalukujen hieman tarkempi lukeminen, luvut 3-4.0, hieman lukua 4.1. Illalla itsearvio"

$ws.Range("D46").WrapText = $true
$ws.Range("D46").Value = "Pääsin paremmin kiinni mikä tilanne oli ennen artikkelia, ja sain pintapuolisen ymmärryksen toisesta artikkelin teesistä. Mielestäni kohtuullinen saavutus 4 tunnin työllä."

$ws.Range("G46").Value = 3.5

$ws.Rows.Item(46).RowHeight = 75

# ---- View state: scroll down & select E46 (matches where the author left off) ----
$ws.Range("E46").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
